# NYPD CompStat weekly sheet refresh: new crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 30   Number  5" -> "Volume 30   Number  6"
$ws.Range("A8").Value = "Volume 30   Number  6"

# "Report Covering the Week  1/30/2023  Through  2/5/2023"
#   -> "Report Covering the Week  2/6/2023  Through  2/12/2023"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Manhattan North weekly crime table (rows 14-30) ---------------------
# Row 14: Murder (label unchanged)
$ws.Range("C14").Value = 2
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -16.666666666666
$ws.Range("I14").Value = 7
$ws.Range("J14").Value = 8
$ws.Range("K14").Value = -12.5
$ws.Range("L14").Value = 75
$ws.Range("M14").Value = -22.222222222222
$ws.Range("N14").Value = -85.416666666666

# Row 15: Rape
$ws.Range("A15").Value = "Rape"
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 11
$ws.Range("H15").Value = -27.272727272727
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 16
$ws.Range("K15").Value = -18.75
$ws.Range("L15").Value = 8.333333333333
$ws.Range("M15").Value = -31.578947368421
$ws.Range("N15").Value = -75.471698113207

# Row 16: Robbery
$ws.Range("A16").Value = "Robbery"
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 48
$ws.Range("E16").Value = -60.416666666666
$ws.Range("F16").Value = 137
$ws.Range("G16").Value = 146
$ws.Range("H16").Value = -6.164383561643
$ws.Range("I16").Value = 214
$ws.Range("J16").Value = 219
$ws.Range("K16").Value = -2.283105022831
$ws.Range("L16").Value = 27.380952380952
$ws.Range("M16").Value = -22.181818181818
$ws.Range("N16").Value = -80.703336339044

# Row 17: Fel. Assault
$ws.Range("A17").Value = "Fel. Assault"
$ws.Range("C17").Value = 54
$ws.Range("D17").Value = 55
$ws.Range("E17").Value = -1.818181818181
$ws.Range("F17").Value = 195
$ws.Range("G17").Value = 195
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 305
$ws.Range("J17").Value = 288
$ws.Range("K17").Value = 5.902777777777
$ws.Range("L17").Value = 17.760617760617
$ws.Range("M17").Value = 50.990099009901
$ws.Range("N17").Value = -44.545454545454

# Row 18: Burglary
$ws.Range("A18").Value = "Burglary"
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = 34
$ws.Range("E18").Value = -11.764705882352
$ws.Range("F18").Value = 104
$ws.Range("G18").Value = 144
$ws.Range("H18").Value = -27.777777777777
$ws.Range("I18").Value = 180
$ws.Range("J18").Value = 240
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -15.492957746478
$ws.Range("M18").Value = -10.447761194029
$ws.Range("N18").Value = -86.384266263237

# Row 19: Gr. Larceny
$ws.Range("A19").Value = "Gr. Larceny"
$ws.Range("C19").Value = 112
$ws.Range("D19").Value = 132
$ws.Range("E19").Value = -15.151515151515
$ws.Range("F19").Value = 434
$ws.Range("G19").Value = 457
$ws.Range("H19").Value = -5.032822757111
$ws.Range("I19").Value = 635
$ws.Range("J19").Value = 673
$ws.Range("K19").Value = -5.646359583952
$ws.Range("L19").Value = 32.845188284518
$ws.Range("M19").Value = 25.494071146245
$ws.Range("N19").Value = -50.851393188854

# Row 20: G.L.A.
$ws.Range("A20").Value = "G.L.A."
$ws.Range("C20").Value = 24
$ws.Range("D20").Value = 38
$ws.Range("E20").Value = -36.842105263157
$ws.Range("F20").Value = 101
$ws.Range("G20").Value = 113
$ws.Range("H20").Value = -10.619469026548
$ws.Range("I20").Value = 145
$ws.Range("J20").Value = 175
$ws.Range("K20").Value = -17.142857142857
$ws.Range("L20").Value = 113.235294117647
$ws.Range("M20").Value = 168.518518518519
$ws.Range("N20").Value = -87.553648068669

# Row 21: TOTAL
$ws.Range("A21").Value = "TOTAL"
$ws.Range("C21").Value = 242
$ws.Range("D21").Value = 309
$ws.Range("E21").Value = -21.682847896440
$ws.Range("F21").Value = 984
$ws.Range("G21").Value = 1072
$ws.Range("H21").Value = -8.208955223880
$ws.Range("I21").Value = 1499
$ws.Range("J21").Value = 1619
$ws.Range("K21").Value = -7.411982705373
$ws.Range("L21").Value = 24.708818635607
$ws.Range("M21").Value = 18.404423380726
$ws.Range("N21").Value = -72.937353312872

# Row 22: Transit
$ws.Range("A22").Value = "Transit"
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 27
$ws.Range("H22").Value = -29.629629629629
$ws.Range("I22").Value = 31
$ws.Range("J22").Value = 43
$ws.Range("K22").Value = -27.906976744186
$ws.Range("L22").Value = 47.619047619047
$ws.Range("M22").Value = -3.125
$ws.Range("N22").Value = "***.*"

# Row 23: Housing
$ws.Range("A23").Value = "Housing"
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 25
$ws.Range("E23").Value = -12
$ws.Range("F23").Value = 90
$ws.Range("G23").Value = 97
$ws.Range("H23").Value = -7.216494845360
$ws.Range("I23").Value = 133
$ws.Range("J23").Value = 138
$ws.Range("K23").Value = -3.623188405797
$ws.Range("L23").Value = 3.100775193798
$ws.Range("M23").Value = 44.565217391304
$ws.Range("N23").Value = "***.*"

# Row 24: Petit Larceny
$ws.Range("A24").Value = "Petit Larceny"
$ws.Range("C24").Value = 269
$ws.Range("D24").Value = 261
$ws.Range("E24").Value = 3.065134099616
$ws.Range("F24").Value = 1083
$ws.Range("G24").Value = 984
$ws.Range("H24").Value = 10.060975609756
$ws.Range("I24").Value = 1617
$ws.Range("J24").Value = 1483
$ws.Range("K24").Value = 9.035738368172
$ws.Range("L24").Value = 32.758620689655
$ws.Range("M24").Value = 58.219178082191
$ws.Range("N24").Value = "***.*"

# Row 25: Misd. Assault
$ws.Range("A25").Value = "Misd. Assault"
$ws.Range("C25").Value = 92
$ws.Range("D25").Value = 86
$ws.Range("E25").Value = 6.976744186046
$ws.Range("F25").Value = 343
$ws.Range("G25").Value = 324
$ws.Range("H25").Value = 5.864197530864
$ws.Range("I25").Value = 489
$ws.Range("J25").Value = 492
$ws.Range("K25").Value = -0.609756097560
$ws.Range("L25").Value = 33.606557377049
$ws.Range("M25").Value = -13.604240282685
$ws.Range("N25").Value = "***.*"

# Row 26: UCR Rape*
$ws.Range("A26").Value = "UCR Rape*"
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -83.333333333333
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 23
$ws.Range("K26").Value = 30.434782608695
$ws.Range("L26").Value = 50
$ws.Range("M26").Value = "***.*"
$ws.Range("N26").Value = "***.*"

# Row 27: Other Sex Crimes
$ws.Range("A27").Value = "Other Sex Crimes"
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 13
$ws.Range("E27").Value = -15.384615384615
$ws.Range("F27").Value = 42
$ws.Range("G27").Value = 47
$ws.Range("H27").Value = -10.638297872340
$ws.Range("I27").Value = 63
$ws.Range("J27").Value = 69
$ws.Range("K27").Value = -8.695652173913
$ws.Range("L27").Value = 23.529411764705
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"

# Row 28: Shooting Vic.
$ws.Range("A28").Value = "Shooting Vic."
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = -31.25
$ws.Range("I28").Value = 19
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -20.833333333333
$ws.Range("L28").Value = 46.153846153846
$ws.Range("M28").Value = 18.75
$ws.Range("N28").Value = -81.904761904761

# Row 29: Shooting Inc.
$ws.Range("A29").Value = "Shooting Inc."
$ws.Range("C29").Value = 3
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 200
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 13
$ws.Range("H29").Value = -23.076923076923
$ws.Range("I29").Value = 18
$ws.Range("J29").Value = 21
$ws.Range("K29").Value = -14.285714285714
$ws.Range("L29").Value = 50
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = -81.818181818181

# Row 30: Hate Crimes
$ws.Range("A30").Value = "Hate Crimes"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 6
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 5
$ws.Range("J30").Value = 8
$ws.Range("K30").Value = -37.5
$ws.Range("L30").Value = -28.571428571428
$ws.Range("M30").Value = "***.*"
$ws.Range("N30").Value = "***.*"
